$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Range Status" sheet: the per-range species counts collapsed to zero and
# the percentage column is dropped (re-ran classify+summarise against the
# changed mapping file, which no longer yields any range-status species).
# ---------------------------------------------------------------------------
$wsRange = $wb.Worksheets.Item("Range Status")

$rangeRows = 2..7
foreach ($r in $rangeRows) {
    $wsRange.Cells.Item($r, 2).Value = 0      # column B -> 0
    $wsRange.Cells.Item($r, 3).Value = $null  # column C -> removed entirely
}

# ---------------------------------------------------------------------------
# "High Priority break-up" sheet: re-summarised numbers, and the "Trend
# Different" / "Range" rows are gone (collapsed into the updated table that
# now only has "Trend New" and "IUCN").
# ---------------------------------------------------------------------------
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")

# Drop the old rows 4 ("Range") and 5 ("IUCN") - delete bottom-up so the
# remaining row indices don't shift under us.
$wsBreakup.Rows.Item(5).Delete()
$wsBreakup.Rows.Item(4).Delete()

# Row 2: "Trend New" gets refreshed counts.
$wsBreakup.Cells.Item(2, 2).Value = 12     # B2 High Species (no.)
$wsBreakup.Cells.Item(2, 3).Value = 21.4   # C2 High Species (perc.)
$wsBreakup.Cells.Item(2, 4).Value = 12     # D2 New High Species (no.)
$wsBreakup.Cells.Item(2, 5).Value = 21.4   # E2 New High Species (perc.)

# Row 3: used to be "Trend Different"; now it's the "IUCN" row with its own
# refreshed counts (and now also has New High Species figures).
$wsBreakup.Cells.Item(3, 1).Value = "IUCN" # A3
$wsBreakup.Cells.Item(3, 2).Value = 44     # B3 High Species (no.)
$wsBreakup.Cells.Item(3, 3).Value = 78.6   # C3 High Species (perc.)
$wsBreakup.Cells.Item(3, 4).Value = 44     # D3 New High Species (no.)
$wsBreakup.Cells.Item(3, 5).Value = 78.6   # E3 New High Species (perc.)
